$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates: volume number and report week dates
$ws.Range("A8").Value = "Volume 30   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/24/2023  Through  4/30/2023"

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 150

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = 17.391304347826
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = 81.818181818181
$ws.Range("M16").Value = -6.976744186046
$ws.Range("N16").Value = -77.186311787072

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 52
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = 40.54054054054
$ws.Range("I17").Value = 202
$ws.Range("J17").Value = 154
$ws.Range("K17").Value = 31.168831168831
$ws.Range("L17").Value = 74.137931034482
$ws.Range("M17").Value = 124.444444444444
$ws.Range("N17").Value = 48.529411764705

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = 35.849056603773
$ws.Range("L18").Value = 26.315789473684
$ws.Range("M18").Value = -41.935483870967
$ws.Range("N18").Value = -91.022443890274

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = 15.625
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 430
$ws.Range("K19").Value = -33.720930232558
$ws.Range("L19").Value = 99.300699300699
$ws.Range("M19").Value = 55.737704918032
$ws.Range("N19").Value = -18.803418803418

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 91
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = 54.237288135593
$ws.Range("L20").Value = 97.826086956521
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = -86.695906432748

# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 7.5
$ws.Range("F21").Value = 192
$ws.Range("G21").Value = 163
$ws.Range("H21").Value = 17.791411042944
$ws.Range("I21").Value = 782
$ws.Range("J21").Value = 805
$ws.Range("K21").Value = -2.857142857142
$ws.Range("L21").Value = 76.126126126126
$ws.Range("M21").Value = 32.76740237691
$ws.Range("N21").Value = -68.894192521877

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 27.272727272727
$ws.Range("M22").Value = 16.666666666666

# Row 24
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 204
$ws.Range("G24").Value = 180
$ws.Range("H24").Value = 13.333333333333
$ws.Range("I24").Value = 924
$ws.Range("J24").Value = 762
$ws.Range("K24").Value = 21.259842519685
$ws.Range("L24").Value = 78.723404255319
$ws.Range("M24").Value = 86.290322580645

# Row 25
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 12.5
$ws.Range("F25").Value = 87
$ws.Range("H25").Value = 45
$ws.Range("I25").Value = 341
$ws.Range("J25").Value = 241
$ws.Range("K25").Value = 41.49377593361
$ws.Range("L25").Value = 54.298642533936
$ws.Range("M25").Value = 72.222222222222

# Row 26
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = -11.764705882352
$ws.Range("L26").Value = -25

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = 28.125
$ws.Range("L27").Value = 24.242424242424

# Row 28
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = -86.363636363636

# Row 29
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -90.47619047619
